$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 18045.28
$ws.Range("B6").Value = 24032.91
$ws.Range("B7").Value = 27515.9
$ws.Range("B13").Value = 15067.65
$ws.Range("B15").Value = 2663.05
$ws.Range("B16").Value = 10624.8
